# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E23) is re-sorted from descending
# (2307 .. 2212) to ascending (2212 .. 2307) chronological order, and the
# "Valor Mora" values in column F follow the same row (F16 <-> F23 swap
# their values so the 2212 row keeps 40000 and the 2307 row keeps 16000).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New ascending order of "Periodo Mora" text for rows 16..23
$periodos = @("2212", "2301", "2302", "2303", "2304", "2305", "2306", "2307")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
}

# "Valor Mora" (column F) swap: row 16 becomes 40000, row 23 becomes 16000
$ws.Cells.Item(16, 6).Value = 40000
$ws.Cells.Item(23, 6).Value = 16000
